$p = $ppt.ActivePresentation

# Slide 4, "TextBox 1" holds the Java code sample that starts with
# "ArrayList<Integer> list = ...". The first "ArrayList" occurrence
# (naming the declared variable's type) is being shortened to "List".
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

$fullText = $tr.Text
$idx = $fullText.IndexOf("ArrayList")
if ($idx -ge 0) {
    $run = $tr.Characters($idx + 1, [string]"ArrayList".Length)
    $run.Text = "List"
}
